$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.936.79'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.637.51'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.86%  '
$ws.Range('D5').Value = '214.79'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = '1.864.59'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.668.52'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '4.24'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '0.544'
$ws.Range('E15').Value = '  -1.58%  '
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').Value = '62.49'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').Value = '25.945.33'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = '193.35'
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('E21').Value = '  -1.64%  '
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('E23').Value = '  -1.44%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').Value = '143.98'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  +0.96%  '
$ws.Range('E27').Value = '  +2.59%  '
$ws.Range('D28').Value = '6.85'
$ws.Range('E28').Value = '  -0.46%  '
$ws.Range('D29').Value = '15.46'
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0500'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').Value = '3.22'
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').Value = '1.136.25'
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '99.26'
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.797'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.41'
$ws.Range('E44').Value = '  -2.56%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.773.29'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0115'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '56.44'
$ws.Range('E47').Value = '  +1.68%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0529'
$ws.Range('E48').Value = '  +2.44%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.45'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.415'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.65'
$ws.Range('E51').Value = '  +0.90%  '
